$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# 1. Text / value updates (schedule content rewrite)
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Stage 1:  Research & Exploratory Measure"
$ws.Range("A3").Value = "4 & 5"
$ws.Range("B4").Value = "Stage 2: Prototype"
$ws.Range("B6").Value = "Stage 3:Basic Skeleton of the App"
$ws.Range("B8").Value = "Stage 4: Minimum Viable Product"
$ws.Range("B11").Value = "Stage 5: Integrated Testing & Stabilisation"

# ---------------------------------------------------------------------------
# 2. Merge layout changes: extend the Stage-4 merge down into row 10 and add
#    a new merged helper column (F) alongside it.
# ---------------------------------------------------------------------------
$ws.Range("B8:B9").UnMerge()
$ws.Range("B8:B10").Merge()
$ws.Range("F8:F10").Merge()

# ---------------------------------------------------------------------------
# 3. Row height tweak (row 3 shrinks a touch once the header text is final)
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 12.5

Write-Output "edit applied"
